$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicaciones")

# Correccion de ortografia: "Descripcion" -> "Descripción"
$ws.Range("C3").Value = "Descripción"

# Update selection to match the saved view state
$ws.Range("C3").Select()
